# Updates "想去人数" (interest count) values in column F
# on sheets "展览" and "全部类型", per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new value }
$updates = @{
    "展览" = @{
        2  = 15166
        3  = 19517
        5  = 162
        14 = 208
        16 = 74
        17 = 1516
        22 = 8177
        24 = 42
        28 = 15
        31 = 6526
        32 = 131
        36 = 306
        37 = 5563
        38 = 1016
    }
    "全部类型" = @{
        2  = 15166
        3  = 19517
        5  = 162
        14 = 208
        16 = 74
        17 = 1516
        23 = 8177
        25 = 42
        29 = 15
        34 = 6526
        35 = 131
        39 = 306
        40 = 5563
        41 = 1016
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Range("F$row").Value = $rowMap[$row]
    }
}
